# Intégration des éléments concernant la structure issus du GT Structure
# - Refresh the generation Date + concept Count on the Metadata sheet.
# - Replace the flat "Entités" concept list with the new "Structure" concept
#   list on the Concepts sheet (Code/Display only - Level stays "1" for every
#   row and is left untouched), and blank out the now-unused Definition
#   column while keeping its existing cell formatting.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2026-01-28T10:29:57+00:00"

# "14" looks numeric, so a plain .Value assignment would turn the cell into
# a number and drop its text type. Write it as text, then restore the
# original (untouched) cell formatting from a neighbouring cell so no new
# style gets introduced.
$meta.Range("B23").NumberFormat = "@"
$meta.Range("B23").Value = "14"
$meta.Range("B22").Copy() | Out-Null
$meta.Range("B23").PasteSpecial(-4122) | Out-Null

$concepts = $wb.Worksheets.Item("Concepts")

# New Code / Display pairs (Level column A is unchanged, always "1").
$codes = @(
    "GHT", "Groupement hospitalier de territoire",
    "LEGAL-ENTITY", "Entité légale",
    "GEOGRAPHICAL-ENTITY", "Entité géographique",
    "GROUP", "Groupe privé/hospitalier",
    "STRUCT-INTERNE", "Structure interne",
    "SECTEUR", "Secteur",
    "DEPARTEMENT", "Département",
    "SERVICE", "Service",
    "UM", "Unité médicale",
    "UAC", "Unité d'activité",
    "POLE", "Pôle",
    "CENTRE-RESP", "Centre de responsabilité",
    "CENTRE-ACTIVITE", "Centre d'activité",
    "UF", "Unité fonctionnelle"
)

$rowCount = $codes.Length / 2

for ($i = 0; $i -lt $rowCount; $i++) {
    $row = $i + 2
    $code = $codes[$i * 2]
    $display = $codes[$i * 2 + 1]
    $concepts.Range("B" + $row).Value = $code
    $concepts.Range("C" + $row).Value = $display
}

# Definition column is no longer populated for any remaining row.
$concepts.Range("D2:D18").ClearContents()

# Only 14 data rows remain (2..15); drop the trailing 3 rows entirely.
$concepts.Rows("16:18").Delete()
